$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the leadership roster names (Column A) for Winter/Spring 2022 EC
$ws.Range("A2").Value = "Grace Hsiang"
$ws.Range("A3").Value = "Aayush Shah"
$ws.Range("A4").Value = "Grace Till"
$ws.Range("A5").Value = "Dylan Tanzil"
$ws.Range("A6").Value = "Tiffany Than"
$ws.Range("A7").Value = "Alex Pham"
$ws.Range("A8").Value = "Megha Bhargava"
$ws.Range("A9").Value = "Tommy Truong"
$ws.Range("A10").Value = "Kevin Cao"
$ws.Range("A11").Value = "Sean Devine"

# Update the active selection state as left by the author (multi-area
# selection: C11 selected first, then D13 ctrl-selected as the active cell)
$ws.Range("C11").Select()
$ws.Range("D13").Select()
